$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:B46").ClearContents() | Out-Null

$ws.Range("A1").Value = "Cluster name"
$ws.Range("B1").Value = "Active cases"

$ws.Range("A2").Value = "3035 Campbell Place Aged Care Glen Waverley"
$ws.Range("B2").Value = 11
$ws.Range("A3").Value = "3528 Ottoman Village Aged Care Broadmeadows"
$ws.Range("B3").Value = 14
$ws.Range("A4").Value = "3622 Olivet Care Aged Care Services Ringwood"
$ws.Range("B4").Value = 13
$ws.Range("A5").Value = "3652 Regis Aged Care Dandenong North"
$ws.Range("B5").Value = 18
$ws.Range("A6").Value = "3824 Estia Health South Morang"
$ws.Range("B6").Value = 27
$ws.Range("A7").Value = "3961 Heritage Care Water Gardens Aged Care Facility Sydenham"
$ws.Range("B7").Value = 15
$ws.Range("A8").Value = "Aintree Primary School Aintree"
$ws.Range("B8").Value = 14
$ws.Range("A9").Value = "Australian Meat Group Abattoir Dandenong South"
$ws.Range("B9").Value = 18
$ws.Range("A10").Value = "Berwick Lodge Primary School Berwick"
$ws.Range("B10").Value = 21
$ws.Range("A11").Value = "Dandenong North Primary School Dandenong"
$ws.Range("B11").Value = 10
$ws.Range("A12").Value = "Elements Childcare Warralily Armstrong Creek"
$ws.Range("B12").Value = 11
$ws.Range("A13").Value = "Hamlyn Views School Hamlyn Heights"
$ws.Range("B13").Value = 12
$ws.Range("A14").Value = "KingKids Early Learning Centre and Kindergarten Hallam"
$ws.Range("B14").Value = 11
$ws.Range("A15").Value = "Lilydale Motor Inn Lilydale"
$ws.Range("B15").Value = 11
$ws.Range("A16").Value = "Lowanna College Newborough"
$ws.Range("B16").Value = 30
$ws.Range("A17").Value = "M.C. Herd Corio"
$ws.Range("B17").Value = 10
$ws.Range("A18").Value = "Melbourne Alternative Place of Detention Carlton"
$ws.Range("B18").Value = 17
$ws.Range("A19").Value = "Morwell Park Primary School Morwell"
$ws.Range("B19").Value = 10
$ws.Range("A20").Value = "Northern Bay College Goldsworthy 9-12 Campus Corio"
$ws.Range("B20").Value = 16
$ws.Range("A21").Value = "Northern Bay College Wexford Campus Corio"
$ws.Range("B21").Value = 48
$ws.Range("A22").Value = "Saint Augustines Primary School Wodonga"
$ws.Range("B22").Value = 10
$ws.Range("A23").Value = "Saint Monica's Primary School Wodonga"
$ws.Range("B23").Value = 12
$ws.Range("A24").Value = "St Josephs Catholic Primary School Warragul"
$ws.Range("B24").Value = 10
$ws.Range("A25").Value = "St Mary's Primary School Swan Hill"
$ws.Range("B25").Value = 21
$ws.Range("A26").Value = "St Mary's Primary School Swan Hill"
$ws.Range("B26").Value = 24
$ws.Range("A27").Value = "St Thereses Primary School Kennington"
$ws.Range("B27").Value = 14
$ws.Range("A28").Value = "St Vincents Hospital Emergency Department Melbourne"
$ws.Range("B28").Value = 15
$ws.Range("A29").Value = "St. Brendans Catholic Primary School Lakes Entrance"
$ws.Range("B29").Value = 10
$ws.Range("A30").Value = "The Royal Children's Hospital Parkville"
$ws.Range("B30").Value = 10
$ws.Range("A31").Value = "Werribee Mercy Hospital Emergency Department"
$ws.Range("B31").Value = 37
$ws.Range("A32").Value = "Western Health Sunshine Hospital Emergency Department St Albans"
$ws.Range("B32").Value = 12
$ws.Range("A33").Value = "Wodonga Cemetery Wodonga Outbreak"
$ws.Range("B33").Value = 38
$ws.Range("A34").Value = "Wodonga Primary School Wodonga"
$ws.Range("B34").Value = 24
$ws.Range("A35").Value = "Wodonga Senior Secondary College Wodonga"
$ws.Range("B35").Value = 24
$ws.Range("A36").Value = "Wodonga South Primary School Wodonga"
$ws.Range("B36").Value = 35
$ws.Range("A37").Value = "Woodend Primary School Woodend"
$ws.Range("B37").Value = 20
$ws.Range("A38").Value = "Yallourn Power Station Yallourn"
$ws.Range("B38").Value = 10
$ws.Range("A39").Value = "Yooralla Disability Residential Care Alfrieda Street St Albans"
$ws.Range("B39").Value = 12
